$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26: I26 "sd" -> "b", J26 "Statement-non-opinion" -> "Acknowledge (Backchannel)"
$ws.Range("I26").Value = "b"
$ws.Range("J26").Value = "Acknowledge (Backchannel)"

# Row 27: I27 "sd" -> "sv", J27 "Statement-non-opinion" -> "Statement-opinion"
$ws.Range("I27").Value = "sv"
$ws.Range("J27").Value = "Statement-opinion"

# Row 43: I43 "sd" -> "aa", J43 "Statement-non-opinion" -> "Agree/Accept"
$ws.Range("I43").Value = "aa"
$ws.Range("J43").Value = "Agree/Accept"

# Row 44: I44 "sd" -> "aa", J44 "Statement-non-opinion" -> "Agree/Accept"
$ws.Range("I44").Value = "aa"
$ws.Range("J44").Value = "Agree/Accept"

# Row 57: I57 "%" -> "aa", J57 "Uninterpretable" -> "Agree/Accept"
$ws.Range("I57").Value = "aa"
$ws.Range("J57").Value = "Agree/Accept"

$wb.Save()
